$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    16.86991607391245,
    15.89584955866815,
    15.2699775710849,
    15.008197319934,
    14.96433081551589,
    15.26647399323134,
    16.53996406344772,
    19.0027458068253,
    20.67494806633232,
    21.39172256362241,
    21.65686569030329,
    21.60004134736741,
    21.4136618050453,
    21.29868154950794,
    20.62722412089977,
    20.20408069597326,
    19.95656407809801,
    19.87204792380568,
    20.24955283636154,
    21.46857628470571,
    22.22866616901554,
    21.82633154458858,
    20.22900810905285,
    18.34778573295695
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
